$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B154").Value = 53925
$ws.Range("E154").Value = 79.37
$ws.Range("F154").Value = 1
$ws.Range("G154").Value = 66.44
$ws.Range("B155").Value = 64350
$ws.Range("E155").Value = 70.63
$ws.Range("F155").Value = 101
$ws.Range("G155").Value = 6710.44
$ws.Range("B156").Value = 57756
$ws.Range("F156").Value = -100
$ws.Range("G156").Value = -6644
$ws.Range("B256").Value = 48719
$ws.Range("E256").Value = 353.35
$ws.Range("F256").Value = -81
$ws.Range("G256").Value = -23955.75
$ws.Range("B257").Value = 64979
$ws.Range("E257").Value = 314.41
$ws.Range("F257").Value = 82
$ws.Range("G257").Value = 24251.5
$ws.Range("B271").Value = 64973
$ws.Range("E271").Value = 35.4
$ws.Range("F271").Value = 150
$ws.Range("G271").Value = 4995
$ws.Range("B272").Value = 48706
$ws.Range("E272").Value = 39.8
$ws.Range("F272").Value = -144
$ws.Range("G272").Value = -4795.2
$ws.Range("B305").Value = 62997
$ws.Range("F305").Value = 72
$ws.Range("G305").Value = 22020.48
$ws.Range("B306").Value = 57854
$ws.Range("F306").Value = 2
$ws.Range("G306").Value = 611.6799999999999
$ws.Range("B309").Value = 61610
$ws.Range("E309").Value = 122.71
$ws.Range("F309").Value = -58
$ws.Range("G309").Value = -5957.18
$ws.Range("B310").Value = 63565
$ws.Range("E310").Value = 109.19
$ws.Range("F310").Value = 60
$ws.Range("G310").Value = 6162.6
$ws.Range("B342").Value = 63531
$ws.Range("E342").Value = 152.53
$ws.Range("F342").Value = 80
$ws.Range("G342").Value = 11478.4
$ws.Range("B343").Value = 57802
$ws.Range("E343").Value = 162.71
$ws.Range("F343").Value = -79
$ws.Range("G343").Value = -11334.92
$ws.Range("B347").Value = 55356
$ws.Range("E347").Value = 54.04
$ws.Range("F347").Value = -158
$ws.Range("G347").Value = -7527.12
$ws.Range("B348").Value = 63510
$ws.Range("E348").Value = 50.66
$ws.Range("F348").Value = 167
$ws.Range("G348").Value = 7955.88
$ws.Range("B367").Value = 63563
$ws.Range("E367").Value = 119.04
$ws.Range("F367").Value = 15
$ws.Range("G367").Value = 1679.4
$ws.Range("B368").Value = 61605
$ws.Range("E368").Value = 133.78
$ws.Range("F368").Value = -13
$ws.Range("G368").Value = -1455.48
$ws.Range("B371").Value = 61608
$ws.Range("E371").Value = 154.12
$ws.Range("F371").Value = -56
$ws.Range("G371").Value = -7224.56
$ws.Range("B372").Value = 63564
$ws.Range("E372").Value = 137.16
$ws.Range("F372").Value = 57
$ws.Range("G372").Value = 7353.57
$ws.Range("B381").Value = 62865
$ws.Range("F381").Value = 151
$ws.Range("G381").Value = 12051.31
$ws.Range("B382").Value = 57817
$ws.Range("F382").Value = 3
$ws.Range("G382").Value = 239.43
$ws.Range("B392").Value = 62933
$ws.Range("F392").Value = 146
$ws.Range("G392").Value = 8632.98
$ws.Range("B393").Value = 57835
$ws.Range("F393").Value = 1
$ws.Range("G393").Value = 59.13
$ws.Range("B411").Value = 57856
$ws.Range("F411").Value = 2
$ws.Range("G411").Value = 342.66
$ws.Range("B412").Value = 63007
$ws.Range("F412").Value = 984
$ws.Range("G412").Value = 168588.72
$ws.Range("B413").Value = 57857
$ws.Range("F413").Value = 3
$ws.Range("G413").Value = 453.51
$ws.Range("B414").Value = 63008
$ws.Range("F414").Value = 504
$ws.Range("G414").Value = 76189.67999999999
$ws.Range("B449").Value = 31930
$ws.Range("E449").Value = 26.8
$ws.Range("F449").Value = -62
$ws.Range("G449").Value = -1390.04
$ws.Range("B450").Value = 63681
$ws.Range("E450").Value = 23.84
$ws.Range("F450").Value = 65
$ws.Range("G450").Value = 1457.3
$ws.Range("B571").Value = 53757
$ws.Range("E571").Value = 16.08
$ws.Range("F571").Value = -159
$ws.Range("G571").Value = -2138.55
$ws.Range("B572").Value = 65069
$ws.Range("E572").Value = 14.3
$ws.Range("F572").Value = 172
$ws.Range("G572").Value = 2313.4
$ws.Range("B575").Value = 65066
$ws.Range("E575").Value = 13.61
$ws.Range("F575").Value = 313
$ws.Range("G575").Value = 4009.53
$ws.Range("B576").Value = 53263
$ws.Range("E576").Value = 15.29
$ws.Range("F576").Value = -309
$ws.Range("G576").Value = -3958.29
$ws.Range("B578").Value = 45695
$ws.Range("E578").Value = 23.58
$ws.Range("F578").Value = -36
$ws.Range("G578").Value = -710.28
$ws.Range("B579").Value = 64915
$ws.Range("E579").Value = 20.98
$ws.Range("F579").Value = 40
$ws.Range("G579").Value = 789.2
$ws.Range("B596").Value = 65067
$ws.Range("E596").Value = 15.65
$ws.Range("F596").Value = 338
$ws.Range("G596").Value = 4978.74
$ws.Range("B597").Value = 53595
$ws.Range("E597").Value = 17.61
$ws.Range("F597").Value = -335
$ws.Range("G597").Value = -4934.55
$ws.Range("B701").Value = 60025
$ws.Range("E701").Value = 37.22
$ws.Range("F701").Value = -98
$ws.Range("G701").Value = -3217.34
$ws.Range("B702").Value = 64833
$ws.Range("E702").Value = 34.9
$ws.Range("F702").Value = 99
$ws.Range("G702").Value = 3250.17
$ws.Range("B712").Value = 60022
$ws.Range("E712").Value = 37.22
$ws.Range("F712").Value = -113
$ws.Range("G712").Value = -3709.79
$ws.Range("B713").Value = 64830
$ws.Range("E713").Value = 34.9
$ws.Range("F713").Value = 117
$ws.Range("G713").Value = 3841.11
$ws.Range("B864").Value = 54751
$ws.Range("E864").Value = 46.34
$ws.Range("F864").Value = -19
$ws.Range("G864").Value = -776.53
$ws.Range("B865").Value = 65079
$ws.Range("E865").Value = 43.44
$ws.Range("F865").Value = 21
$ws.Range("G865").Value = 858.27

Write-Output "Applied 161 cell updates"
